$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B8").Value = "Not enough time alloted to complete the task"
$ws.Range("C8").Value = "Lack of Time"
$ws.Range("D8").Value = "Extend time allocation"
$ws.Range("E8").Value = "To fully complete the task"
$ws.Range("F8").Value = "Likely"
$ws.Range("G8").Value = "Major"
$ws.Range("H8").Value = 12

$ws.Range("H8").Interior.ThemeColor = 8
